$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 129.14285
$ws.Range("I6").Value = 145.75
$ws.Range("K6").Value = 437.25
$ws.Range("M6").Value = -325.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2224.25
$ws.Range("I31").Value = 2224.25
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6672.75
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -6442.75
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 271.8
$ws.Range("I39").Value = 245.88235
$ws.Range("J39").Value = 418.66666
$ws.Range("K39").Value = 737.64705
$ws.Range("L39").Value = 1255.99998
$ws.Range("M39").Value = -441.64705
$ws.Range("N39").Value = -1847.99998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 9950
$ws.Range("I47").Value = 9950
$ws.Range("K47").Value = 9950
$ws.Range("M47").Value = -8978

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 28403.25
$ws.Range("J69").Value = 15629.267
$ws.Range("L69").Value = 46887.801
$ws.Range("N69").Value = -48635.801

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 28403.25
$ws.Range("J72").Value = 15629.267
$ws.Range("L72").Value = 140663.403
$ws.Range("N72").Value = -149399.403

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3624.1904
$ws.Range("J138").Value = 3831.5334
$ws.Range("L138").Value = 11494.6002
$ws.Range("N138").Value = -21774.6002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1001
$ws.Range("I5").Value = 1001
$ws.Range("K5").Value = 1001
$ws.Range("M5").Value = -889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1801.8948
$ws.Range("I45").Value = 1216.8
$ws.Range("J45").Value = 3996
$ws.Range("K45").Value = 1216.8
$ws.Range("L45").Value = 3996
$ws.Range("M45").Value = -839.8
$ws.Range("N45").Value = -4750

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 30018.637
$ws.Range("I122").Value = 36022.89
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 108068.67
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -105618.67
$ws.Range("N122").Value = -13898.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1835
$ws.Range("I132").Value = 1711.25
$ws.Range("K132").Value = 5133.75
$ws.Range("M132").Value = -2603.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1001
$ws.Range("I4").Value = 1001
$ws.Range("K4").Value = 1001
$ws.Range("M4").Value = -886

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 14000
$ws.Range("I29").Value = 10000
$ws.Range("J29").Value = 18000
$ws.Range("K29").Value = 10000
$ws.Range("L29").Value = 18000
$ws.Range("M29").Value = -9711
$ws.Range("N29").Value = -18578

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6667
$ws.Range("I86").Value = 6667
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 6667
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -5544
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 6667
$ws.Range("I89").Value = 6667
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 33335
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -27719
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1879.5714
$ws.Range("I94").Value = 1901.4445
$ws.Range("J94").Value = 1748.3334
$ws.Range("K94").Value = 1901.4445
$ws.Range("L94").Value = 1748.3334
$ws.Range("M94").Value = -1450.4445
$ws.Range("N94").Value = -2650.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4711.25
$ws.Range("I107").Value = 4398.2
$ws.Range("K107").Value = 4398.2
$ws.Range("M107").Value = -2478.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4496.125
$ws.Range("I134").Value = 4496.125
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13488.375
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -10953.375
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1884.6666
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1884.6666
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1884.6666
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2584.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1075.4762
$ws.Range("I31").Value = 1053.4
$ws.Range("J31").Value = 1082.375
$ws.Range("K31").Value = 1053.4
$ws.Range("L31").Value = 1082.375
$ws.Range("M31").Value = -758.4000000000001
$ws.Range("N31").Value = -1672.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 10430
$ws.Range("I32").Value = 9906.666999999999
$ws.Range("K32").Value = 9906.666999999999
$ws.Range("M32").Value = -9590.666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1075.4762
$ws.Range("I34").Value = 1053.4
$ws.Range("J34").Value = 1082.375
$ws.Range("K34").Value = 1053.4
$ws.Range("L34").Value = 1082.375
$ws.Range("M34").Value = -851.4000000000001
$ws.Range("N34").Value = -1486.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 50079
$ws.Range("I36").Value = 50079
$ws.Range("K36").Value = 50079
$ws.Range("M36").Value = -49691

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 50079
$ws.Range("I40").Value = 50079
$ws.Range("K40").Value = 50079
$ws.Range("M40").Value = -49919

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 38333
$ws.Range("I42").Value = 38333
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 38333
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -37740
$ws.Range("N42").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4681.125
$ws.Range("I86").Value = 4625.1665
$ws.Range("K86").Value = 4625.1665
$ws.Range("M86").Value = -3502.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 4681.125
$ws.Range("I89").Value = 4625.1665
$ws.Range("K89").Value = 23125.8325
$ws.Range("M89").Value = -17509.8325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 41483
$ws.Range("J92").Value = 41483
$ws.Range("L92").Value = 41483
$ws.Range("N92").Value = -46475

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7229
$ws.Range("I132").Value = 7229
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 21687
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -19157
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 3375.4
$ws.Range("I11").Value = 3375.4
$ws.Range("K11").Value = 10126.2
$ws.Range("M11").Value = -9986.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 20330
$ws.Range("I63").Value = 20330
$ws.Range("K63").Value = 60990
$ws.Range("M63").Value = -60241

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 20330
$ws.Range("I66").Value = 20330
$ws.Range("K66").Value = 182970
$ws.Range("M66").Value = -179226

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 236
$ws.Range("J86").Value = 260
$ws.Range("L86").Value = 780
$ws.Range("N86").Value = -3152

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 236
$ws.Range("J89").Value = 260
$ws.Range("L89").Value = 2340
$ws.Range("N89").Value = -14196

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 527805.7
$ws.Range("I131").Value = 887
$ws.Range("J131").Value = 1002032.5
$ws.Range("K131").Value = 2661
$ws.Range("L131").Value = 3006097.5
$ws.Range("M131").Value = 2379
$ws.Range("N131").Value = -3016177.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 20033332
$ws.Range("J33").Value = 20033332
$ws.Range("L33").Value = 20033332
$ws.Range("N33").Value = -20033836

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51872

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 35999.5
$ws.Range("J101").Value = 35999.5
$ws.Range("L101").Value = 35999.5
$ws.Range("N101").Value = -42489.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3299.2
$ws.Range("I102").Value = 2996
$ws.Range("K102").Value = 2996
$ws.Range("M102").Value = -1374

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5255.5
$ws.Range("I126").Value = 4231
$ws.Range("J126").Value = 5511.625
$ws.Range("K126").Value = 12693
$ws.Range("L126").Value = 16534.875
$ws.Range("M126").Value = -10223
$ws.Range("N126").Value = -21474.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6445.3335
$ws.Range("I7").Value = 5241.3335
$ws.Range("J7").Value = 7649.3335
$ws.Range("K7").Value = 5241.3335
$ws.Range("L7").Value = 7649.3335
$ws.Range("M7").Value = -5129.3335
$ws.Range("N7").Value = -7873.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1249.25
$ws.Range("I61").Value = 999.5
$ws.Range("J61").Value = 1499
$ws.Range("K61").Value = 999.5
$ws.Range("L61").Value = 1499
$ws.Range("M61").Value = -797.5
$ws.Range("N61").Value = -1903

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1249.25
$ws.Range("I113").Value = 999.5
$ws.Range("J113").Value = 1499
$ws.Range("K113").Value = 999.5
$ws.Range("L113").Value = 1499
$ws.Range("M113").Value = 1170.5
$ws.Range("N113").Value = -5839

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6445.3335
$ws.Range("I126").Value = 5241.3335
$ws.Range("J126").Value = 7649.3335
$ws.Range("K126").Value = 15724.0005
$ws.Range("L126").Value = 22948.0005
$ws.Range("M126").Value = -13254.0005
$ws.Range("N126").Value = -27888.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 801.3333
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = 4
$ws.Range("M14").Value = 164

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1996.4546
$ws.Range("I107").Value = 1884.5
$ws.Range("J107").Value = 2295
$ws.Range("K107").Value = 5653.5
$ws.Range("L107").Value = 6885
$ws.Range("M107").Value = -3733.5
$ws.Range("N107").Value = -10725
